# The deck currently has the "Integral" (Red Violet) design applied to the
# slide master/theme, plus a second, unused "Office Theme" theme part left
# over in the package. The edit re-applies the stock "Office Theme" colour
# palette to the presentation's (single) live theme -- i.e. it's the
# Design > Themes > "Office Theme" click in the ribbon -- swapping the
# deck's visible colour scheme from Red Violet back to the default Office
# palette. Font scheme / format scheme are already identical between the
# two themes, so only the 12 theme colour slots need to change.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Helper: PowerPoint's ColorFormat.RGB (and the classic VBA RGB() function)
# packs components as 0x00BBGGRR, i.e. R + G*256 + B*65536.
function ToRgb([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in ThemeColorScheme.Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeThemeColors = @(
    (ToRgb 0x00 0x00 0x00),  #  1 dk1      000000
    (ToRgb 0xFF 0xFF 0xFF),  #  2 lt1      FFFFFF
    (ToRgb 0x44 0x54 0x6A),  #  3 dk2      44546A
    (ToRgb 0xE7 0xE6 0xE6),  #  4 lt2      E7E6E6
    (ToRgb 0x5B 0x9B 0xD5),  #  5 accent1  5B9BD5
    (ToRgb 0xED 0x7D 0x31),  #  6 accent2  ED7D31
    (ToRgb 0xA5 0xA5 0xA5),  #  7 accent3  A5A5A5
    (ToRgb 0xFF 0xC0 0x00),  #  8 accent4  FFC000
    (ToRgb 0x44 0x72 0xC4),  #  9 accent5  4472C4
    (ToRgb 0x70 0xAD 0x47),  # 10 accent6  70AD47
    (ToRgb 0x05 0x63 0xC1),  # 11 hlink    0563C1
    (ToRgb 0x95 0x4F 0x72)   # 12 folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
